$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 becomes a new "test log" entry, mirroring the structure/format of row 11
# (31/08/2021, RA, 1, ok) but recording a new sign-in test run.
$ws.Range("B11:F11").Copy()
$ws.Range("B12:F12").PasteSpecial()

$ws.Range("B12").Value = "902/9/2021"
$ws.Range("B12").NumberFormat = "mm-dd-yy"
$ws.Range("C12").Value = "RA"
$ws.Range("D12").Value = 1
$ws.Range("E12").ClearContents()
$ws.Range("F12").Value = "ok"

# Drop the leftover empty placeholder cells in column D (no value, no formatting impact)
$ws.Range("D3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("D27").ClearContents()
$ws.Range("D28").ClearContents()
$ws.Range("D29").ClearContents()
$ws.Range("D32").ClearContents()
$ws.Range("D33").ClearContents()
$ws.Range("D34").ClearContents()
$ws.Range("D35").ClearContents()
$ws.Range("D36").ClearContents()
$ws.Range("D37").ClearContents()
